$wb = $excel.ActiveWorkbook

# Update "Latest HO Xliff Generate Date" on the Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-12 03:29:55"

# Update "Correspond Handoff Datetime" and "Correspond Handback DateTime" on zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-12 03:29:50"
$wsZhCn.Range("K2").Value = "2016-08-12 03:30:17"

# Update "Correspond Handoff Datetime" and "Correspond Handback DateTime" on de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-12 03:29:55"
$wsDeDe.Range("K2").Value = "2016-08-12 03:30:24"
